$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I17").Value = -0.2889366546078636
$ws.Range("J17").Value = 0.1255641108669581
$ws.Range("K17").Value = 0.2492139062772706
$ws.Range("L17").Value = 2.134687693561835
